# "Fixed major bugs in MyCam"
#
# The cam-profile table (rows 7-13) had incorrect Action/Displacement/Degree
# combinations. The corrected workbook collapses the table down to 4 rows
# (rows 7-10), all using Action = RISE, Displacement = 10, Degree = 90, and
# one row per Motion Type (SHM, UARM, CYCLOIDAL, UV). Rows 11-13 are removed
# entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused rows 11-13 completely (content + formatting), so the
# sheet's used range/dimension shrinks back down to A1:D10 and the orphaned
# "DWELL"/"FALL" shared strings naturally fall out of the saved workbook.
$ws.Range("A11:D13").Clear()

# Rewrite the corrected cam-profile rows.
$ws.Range("A7").Value = "RISE"
$ws.Range("B7").Value = "SHM"
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 90

$ws.Range("A8").Value = "RISE"
$ws.Range("B8").Value = "UARM"
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 90

$ws.Range("A9").Value = "RISE"
$ws.Range("B9").Value = "CYCLOIDAL"
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 90

$ws.Range("A10").Value = "RISE"
$ws.Range("B10").Value = "UV"
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = 90

# Update the view: scroll down a little and leave the selection on D11 (just
# below the now-shorter table), matching where the author's cursor ended up.
try {
    $ws.Activate()
    $excel.ActiveWindow.ScrollRow = 4
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scroll-position isn't always settable in every host; selection below
    # is the part that matters for the saved workbook state.
}
$ws.Range("D11").Select()
